# Generate Report for Archive
#
# The localization status of the handed-off file moved on from
# "Ready for handoff" to "In Translation". That status string is shared
# across the Overview sheet (columns zh-cn / de-de) and each per-locale
# sheet's "Status" column, so every cell that shows it needs to be
# refreshed. Updating the text also shortens the column contents, so the
# Status columns are re-sized to fit the new text (mirroring what the
# report generator does when it regenerates this workbook).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Overview sheet: zh-cn / de-de status columns (E, F) on the data row.
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Per-locale sheets: "Status" column (C) on the data row.
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# The shorter status text means the Status columns no longer need to be
# as wide - shrink them to fit the new content, same as the report tool
# does when it rebuilds this file.
$newWidth = 12.5

$overview.Range("E1").EntireColumn.ColumnWidth = $newWidth
$overview.Range("F1").EntireColumn.ColumnWidth = $newWidth
$zhcn.Range("C1").EntireColumn.ColumnWidth = $newWidth
$dede.Range("C1").EntireColumn.ColumnWidth = $newWidth
